$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 31 and 34 (full F:V range swap) ---
$r31 = $ws.Range("F31:V31").Value2
$r34 = $ws.Range("F34:V34").Value2
$ws.Range("F31:V31").Value2 = $r34
$ws.Range("F34:V34").Value2 = $r31

# --- Swap rows 43 and 44 (full F:V range swap) ---
$r43 = $ws.Range("F43:V43").Value2
$r44 = $ws.Range("F44:V44").Value2
$ws.Range("F43:V43").Value2 = $r44
$ws.Range("F44:V44").Value2 = $r43

# --- Append 7 new match rows (55-61) ---

# Row 55 (Indice 54)
$ws.Cells.Item(55, 1).Value2 = 54
$ws.Cells.Item(55, 2).Value2 = "italy"
$ws.Cells.Item(55, 3).Value2 = "serie-b"
$ws.Cells.Item(55, 4).Value2 = "2023-2024"
$ws.Cells.Item(55, 5).Value2 = 45195.76041666666
$ws.Cells.Item(55, 6).Value2 = "Cosenza"
$ws.Cells.Item(55, 7).Value2 = 1
$ws.Cells.Item(55, 8).Value2 = "Cremonese"
$ws.Cells.Item(55, 9).Value2 = 2
$ws.Cells.Item(55, 10).Value2 = 3.16
$ws.Cells.Item(55, 11).Value2 = "24/09/2023 23:12"
$ws.Cells.Item(55, 12).Value2 = 3.6
$ws.Cells.Item(55, 13).Value2 = "26/09/2023 18:10"
$ws.Cells.Item(55, 14).Value2 = 3.3
$ws.Cells.Item(55, 15).Value2 = "24/09/2023 23:12"
$ws.Cells.Item(55, 16).Value2 = 3.28
$ws.Cells.Item(55, 17).Value2 = "26/09/2023 17:55"
$ws.Cells.Item(55, 18).Value2 = 2.39
$ws.Cells.Item(55, 19).Value2 = "24/09/2023 23:12"
$ws.Cells.Item(55, 20).Value2 = 2.25
$ws.Cells.Item(55, 21).Value2 = "26/09/2023 17:52"
$ws.Cells.Item(55, 22).Value2 = "https://www.betexplorer.com/football/italy/serie-b/cosenza-cremonese/8WSDwDI0/"
$ws.Cells.Item(54, 1).Copy()
$ws.Cells.Item(55, 1).PasteSpecial(-4122)
$ws.Cells.Item(54, 5).Copy()
$ws.Cells.Item(55, 5).PasteSpecial(-4122)

# Row 56 (Indice 55)
$ws.Cells.Item(56, 1).Value2 = 55
$ws.Cells.Item(56, 2).Value2 = "italy"
$ws.Cells.Item(56, 3).Value2 = "serie-b"
$ws.Cells.Item(56, 4).Value2 = "2023-2024"
$ws.Cells.Item(56, 5).Value2 = 45195.76041666666
$ws.Cells.Item(56, 6).Value2 = "Lecco"
$ws.Cells.Item(56, 7).Value2 = 1
$ws.Cells.Item(56, 8).Value2 = "FeralpiSalo"
$ws.Cells.Item(56, 9).Value2 = 2
$ws.Cells.Item(56, 10).Value2 = 2.32
$ws.Cells.Item(56, 11).Value2 = "24/09/2023 23:12"
$ws.Cells.Item(56, 12).Value2 = 2.51
$ws.Cells.Item(56, 13).Value2 = "26/09/2023 18:11"
$ws.Cells.Item(56, 14).Value2 = 3.15
$ws.Cells.Item(56, 15).Value2 = "24/09/2023 23:12"
$ws.Cells.Item(56, 16).Value2 = 3.11
$ws.Cells.Item(56, 17).Value2 = "26/09/2023 18:11"
$ws.Cells.Item(56, 18).Value2 = 3.44
$ws.Cells.Item(56, 19).Value2 = "24/09/2023 23:12"
$ws.Cells.Item(56, 20).Value2 = 3.25
$ws.Cells.Item(56, 21).Value2 = "26/09/2023 18:11"
$ws.Cells.Item(56, 22).Value2 = "https://www.betexplorer.com/football/italy/serie-b/lecco-feralpisalo/jeIIxXX6/"
$ws.Cells.Item(54, 1).Copy()
$ws.Cells.Item(56, 1).PasteSpecial(-4122)
$ws.Cells.Item(54, 5).Copy()
$ws.Cells.Item(56, 5).PasteSpecial(-4122)

# Row 57 (Indice 56)
$ws.Cells.Item(57, 1).Value2 = 56
$ws.Cells.Item(57, 2).Value2 = "italy"
$ws.Cells.Item(57, 3).Value2 = "serie-b"
$ws.Cells.Item(57, 4).Value2 = "2023-2024"
$ws.Cells.Item(57, 5).Value2 = 45195.85416666666
$ws.Cells.Item(57, 6).Value2 = "Ascoli"
$ws.Cells.Item(57, 7).Value2 = 2
$ws.Cells.Item(57, 8).Value2 = "Ternana"
$ws.Cells.Item(57, 9).Value2 = 0
$ws.Cells.Item(57, 10).Value2 = 2.54
$ws.Cells.Item(57, 11).Value2 = "23/09/2023 13:13"
$ws.Cells.Item(57, 12).Value2 = 2.41
$ws.Cells.Item(57, 13).Value2 = "26/09/2023 20:26"
$ws.Cells.Item(57, 14).Value2 = 3.3
$ws.Cells.Item(57, 15).Value2 = "23/09/2023 13:13"
$ws.Cells.Item(57, 16).Value2 = 3.15
$ws.Cells.Item(57, 17).Value2 = "26/09/2023 20:26"
$ws.Cells.Item(57, 18).Value2 = 2.99
$ws.Cells.Item(57, 19).Value2 = "23/09/2023 13:13"
$ws.Cells.Item(57, 20).Value2 = 3.39
$ws.Cells.Item(57, 21).Value2 = "26/09/2023 20:26"
$ws.Cells.Item(57, 22).Value2 = "https://www.betexplorer.com/football/italy/serie-b/ascoli-ternana/rPV1tFYs/"
$ws.Cells.Item(54, 1).Copy()
$ws.Cells.Item(57, 1).PasteSpecial(-4122)
$ws.Cells.Item(54, 5).Copy()
$ws.Cells.Item(57, 5).PasteSpecial(-4122)

# Row 58 (Indice 57)
$ws.Cells.Item(58, 1).Value2 = 57
$ws.Cells.Item(58, 2).Value2 = "italy"
$ws.Cells.Item(58, 3).Value2 = "serie-b"
$ws.Cells.Item(58, 4).Value2 = "2023-2024"
$ws.Cells.Item(58, 5).Value2 = 45195.85416666666
$ws.Cells.Item(58, 6).Value2 = "Reggiana"
$ws.Cells.Item(58, 7).Value2 = 0
$ws.Cells.Item(58, 8).Value2 = "Pisa"
$ws.Cells.Item(58, 9).Value2 = 0
$ws.Cells.Item(58, 10).Value2 = 2.76
$ws.Cells.Item(58, 11).Value2 = "25/09/2023 01:12"
$ws.Cells.Item(58, 12).Value2 = 2.33
$ws.Cells.Item(58, 13).Value2 = "26/09/2023 20:29"
$ws.Cells.Item(58, 14).Value2 = 3.1
$ws.Cells.Item(58, 15).Value2 = "25/09/2023 01:12"
$ws.Cells.Item(58, 16).Value2 = 3.08
$ws.Cells.Item(58, 17).Value2 = "26/09/2023 20:29"
$ws.Cells.Item(58, 18).Value2 = 2.88
$ws.Cells.Item(58, 19).Value2 = "25/09/2023 01:12"
$ws.Cells.Item(58, 20).Value2 = 3.64
$ws.Cells.Item(58, 21).Value2 = "26/09/2023 20:29"
$ws.Cells.Item(58, 22).Value2 = "https://www.betexplorer.com/football/italy/serie-b/reggiana-pisa/I1A40zQD/"
$ws.Cells.Item(54, 1).Copy()
$ws.Cells.Item(58, 1).PasteSpecial(-4122)
$ws.Cells.Item(54, 5).Copy()
$ws.Cells.Item(58, 5).PasteSpecial(-4122)

# Row 59 (Indice 58)
$ws.Cells.Item(59, 1).Value2 = 58
$ws.Cells.Item(59, 2).Value2 = "italy"
$ws.Cells.Item(59, 3).Value2 = "serie-b"
$ws.Cells.Item(59, 4).Value2 = "2023-2024"
$ws.Cells.Item(59, 5).Value2 = 45195.85416666666
$ws.Cells.Item(59, 6).Value2 = "Spezia"
$ws.Cells.Item(59, 7).Value2 = 0
$ws.Cells.Item(59, 8).Value2 = "Brescia"
$ws.Cells.Item(59, 9).Value2 = 0
$ws.Cells.Item(59, 10).Value2 = 1.86
$ws.Cells.Item(59, 11).Value2 = "25/09/2023 01:42"
$ws.Cells.Item(59, 12).Value2 = 1.83
$ws.Cells.Item(59, 13).Value2 = "26/09/2023 20:29"
$ws.Cells.Item(59, 14).Value2 = 3.65
$ws.Cells.Item(59, 15).Value2 = "25/09/2023 01:42"
$ws.Cells.Item(59, 16).Value2 = 3.58
$ws.Cells.Item(59, 17).Value2 = "26/09/2023 20:29"
$ws.Cells.Item(59, 18).Value2 = 4.51
$ws.Cells.Item(59, 19).Value2 = "25/09/2023 01:42"
$ws.Cells.Item(59, 20).Value2 = 4.89
$ws.Cells.Item(59, 21).Value2 = "26/09/2023 20:29"
$ws.Cells.Item(59, 22).Value2 = "https://www.betexplorer.com/football/italy/serie-b/spezia-brescia/dUUCKR0a/"
$ws.Cells.Item(54, 1).Copy()
$ws.Cells.Item(59, 1).PasteSpecial(-4122)
$ws.Cells.Item(54, 5).Copy()
$ws.Cells.Item(59, 5).PasteSpecial(-4122)

# Row 60 (Indice 59)
$ws.Cells.Item(60, 1).Value2 = 59
$ws.Cells.Item(60, 2).Value2 = "italy"
$ws.Cells.Item(60, 3).Value2 = "serie-b"
$ws.Cells.Item(60, 4).Value2 = "2023-2024"
$ws.Cells.Item(60, 5).Value2 = 45195.85416666666
$ws.Cells.Item(60, 6).Value2 = "Sudtirol"
$ws.Cells.Item(60, 7).Value2 = 0
$ws.Cells.Item(60, 8).Value2 = "Modena"
$ws.Cells.Item(60, 9).Value2 = 0
$ws.Cells.Item(60, 10).Value2 = 2.79
$ws.Cells.Item(60, 11).Value2 = "23/09/2023 13:13"
$ws.Cells.Item(60, 12).Value2 = 2.86
$ws.Cells.Item(60, 13).Value2 = "26/09/2023 20:29"
$ws.Cells.Item(60, 14).Value2 = 3.1
$ws.Cells.Item(60, 15).Value2 = "23/09/2023 13:13"
$ws.Cells.Item(60, 16).Value2 = 2.8
$ws.Cells.Item(60, 17).Value2 = "26/09/2023 20:29"
$ws.Cells.Item(60, 18).Value2 = 2.8
$ws.Cells.Item(60, 19).Value2 = "23/09/2023 13:13"
$ws.Cells.Item(60, 20).Value2 = 3.11
$ws.Cells.Item(60, 21).Value2 = "26/09/2023 20:29"
$ws.Cells.Item(60, 22).Value2 = "https://www.betexplorer.com/football/italy/serie-b/sudtirol-modena/dh98afuK/"
$ws.Cells.Item(54, 1).Copy()
$ws.Cells.Item(60, 1).PasteSpecial(-4122)
$ws.Cells.Item(54, 5).Copy()
$ws.Cells.Item(60, 5).PasteSpecial(-4122)

# Row 61 (Indice 60)
$ws.Cells.Item(61, 1).Value2 = 60
$ws.Cells.Item(61, 2).Value2 = "italy"
$ws.Cells.Item(61, 3).Value2 = "serie-b"
$ws.Cells.Item(61, 4).Value2 = "2023-2024"
$ws.Cells.Item(61, 5).Value2 = 45195.85416666666
$ws.Cells.Item(61, 6).Value2 = "Venezia"
$ws.Cells.Item(61, 7).Value2 = 1
$ws.Cells.Item(61, 8).Value2 = "Palermo"
$ws.Cells.Item(61, 9).Value2 = 3
$ws.Cells.Item(61, 10).Value2 = 2.24
$ws.Cells.Item(61, 11).Value2 = "23/09/2023 13:13"
$ws.Cells.Item(61, 12).Value2 = 2.49
$ws.Cells.Item(61, 13).Value2 = "26/09/2023 20:28"
$ws.Cells.Item(61, 14).Value2 = 3.46
$ws.Cells.Item(61, 15).Value2 = "23/09/2023 13:13"
$ws.Cells.Item(61, 16).Value2 = 3.25
$ws.Cells.Item(61, 17).Value2 = "26/09/2023 20:28"
$ws.Cells.Item(61, 18).Value2 = 3.37
$ws.Cells.Item(61, 19).Value2 = "23/09/2023 13:13"
$ws.Cells.Item(61, 20).Value2 = 3.14
$ws.Cells.Item(61, 21).Value2 = "26/09/2023 20:28"
$ws.Cells.Item(61, 22).Value2 = "https://www.betexplorer.com/football/italy/serie-b/venezia-palermo/ET3DbEfQ/"
$ws.Cells.Item(54, 1).Copy()
$ws.Cells.Item(61, 1).PasteSpecial(-4122)
$ws.Cells.Item(54, 5).Copy()
$ws.Cells.Item(61, 5).PasteSpecial(-4122)

$excel.CutCopyMode = $false
